$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.133
$ws.Range("E4").Value = 16.545
$ws.Range("B11").Value = 6.731
$ws.Range("B12").Value = 5.492999999999999
$ws.Range("E14").Value = 17.175
$ws.Range("B15").Value = 5.295
$ws.Range("E26").Value = 16.179
$ws.Range("B27").Value = 5.322000000000001
$ws.Range("B28").Value = 6.013
$ws.Range("B31").Value = 6.029000000000001
$ws.Range("E31").Value = 16.196
$ws.Range("B32").Value = 6.745
$ws.Range("E35").Value = 16.545
$ws.Range("B36").Value = 8.507999999999999
$ws.Range("E37").Value = 16.653
$ws.Range("B38").Value = 5.225
$ws.Range("E39").Value = 16.444
$ws.Range("E40").Value = 16.693
$ws.Range("E45").Value = 16.765
$ws.Range("B46").Value = 6.367
$ws.Range("E52").Value = 16.718
$ws.Range("B54").Value = 5.072000000000001
$ws.Range("B55").Value = 4.459
$ws.Range("B56").Value = 4.7
$ws.Range("E57").Value = 16.529
$ws.Range("B67").Value = 4.988999999999999
$ws.Range("B69").Value = 5.002000000000001
$ws.Range("B72").Value = 6.347
$ws.Range("B73").Value = 6.600999999999999
$ws.Range("E81").Value = 16.309
$ws.Range("B83").Value = 5.827
$ws.Range("E83").Value = 16.973
$ws.Range("B86").Value = 5.016
$ws.Range("B91").Value = 6.027
$ws.Range("B93").Value = 5.380999999999999
$ws.Range("B99").Value = 5.217
$ws.Range("E100").Value = 16.568
$ws.Range("E102").Value = 16.383
